$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N, shifting Late/heading/Outstanding
# columns one slot to the right (N->O, O->P, P->Q).
$ws.Columns("N").EntireColumn.Insert()

# The newly inserted column inherits the width of the column to its left (M).
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab and move the selection.
$ws.Activate() | Out-Null
$ws.Range("I19").Select() | Out-Null
